# Update cryptocurrency price (D) and volume-change (E) columns with
# the latest scrape, as text values (matching the source data format).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D column cells whose new text looks like a plain number: pre-format
# the cells as Text so Excel stores them as strings (not doubles), then
# strip the explicit Text format again so no stray cell style is left
# behind (matches the original file, which carries no style on these
# cells).
$numericLooking = @("D4","D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D17","D19","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D37","D40","D41","D42","D43","D44","D45","D47","D49","D50","D51")
foreach ($ref in $numericLooking) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.915.26"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "1.831.44"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "244.22"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "0.6878"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "0.07663"
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("D9").Value = "0.3052"
$ws.Range("E9").Value = "  -2.60%  "
$ws.Range("D10").Value = "23.53"
$ws.Range("E10").Value = "  -4.08%  "
$ws.Range("D11").Value = "0.07799"
$ws.Range("E11").Value = "  -2.21%  "
$ws.Range("D12").Value = "1.840.62"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").Value = "5.078"
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").Value = "90.50"
$ws.Range("E14").Value = "  -3.10%  "
$ws.Range("D15").Value = "0.6781"
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("D16").Value = "6.457"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "0.000008292"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "28.926.85"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").Value = "242.84"
$ws.Range("E19").Value = "  -3.77%  "
$ws.Range("D20").Value = "2.079.15"
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("D21").Value = "12.69"
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "7.443"
$ws.Range("E23").Value = "  -2.37%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "0.1472"
$ws.Range("E25").Value = "  -5.39%  "
$ws.Range("D26").Value = "161.42"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "8.788"
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("D28").Value = "18.18"
$ws.Range("E28").Value = "  -2.86%  "
$ws.Range("D29").Value = "1.541"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("D30").Value = "4.214"
$ws.Range("E30").Value = "  -2.61%  "
$ws.Range("D31").Value = "4.126"
$ws.Range("E31").Value = "  -3.26%  "
$ws.Range("D32").Value = "1.180"
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("D33").Value = "0.05122"
$ws.Range("E33").Value = "  -3.62%  "
$ws.Range("D34").Value = "0.7553"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").Value = "1.833"
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("D37").Value = "2.681"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("D39").Value = "1.230.55"
$ws.Range("E39").Value = "  -3.63%  "
$ws.Range("D40").Value = "2.698"
$ws.Range("E40").Value = "  -1.81%  "
$ws.Range("D41").Value = "0.9199"
$ws.Range("E41").Value = "  +2.79%  "
$ws.Range("D42").Value = "108.49"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").Value = "0.9996"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "5.644"
$ws.Range("E44").Value = "  -7.12%  "
$ws.Range("D45").Value = "0.5169"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").Value = "1.978.62"
$ws.Range("E46").Value = "  -2.42%  "
$ws.Range("D47").Value = "9.504"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("E48").Value = "  -3.91%  "
$ws.Range("D49").Value = "64.28"
$ws.Range("E49").Value = "  -9.96%  "
$ws.Range("D50").Value = "1.737"
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("D51").Value = "0.4196"
$ws.Range("E51").Value = "  -2.69%  "

foreach ($ref in $numericLooking) {
    $ws.Range($ref).ClearFormats()
}
